$d = $word.ActiveDocument

$d.Content.Find.Execute("694×3=", $true, $false, $false, $false, $false, $true, 1, $false, "451×2=", 2) | Out-Null
$d.Content.Find.Execute("102×4=", $true, $false, $false, $false, $false, $true, 1, $false, "464×9=", 2) | Out-Null
$d.Content.Find.Execute("577×8=", $true, $false, $false, $false, $false, $true, 1, $false, "415×8=", 2) | Out-Null
$d.Content.Find.Execute("446×4=", $true, $false, $false, $false, $false, $true, 1, $false, "625×2=", 2) | Out-Null
$d.Content.Find.Execute("276×5=", $true, $false, $false, $false, $false, $true, 1, $false, "479×7=", 2) | Out-Null
$d.Content.Find.Execute("528×9=", $true, $false, $false, $false, $false, $true, 1, $false, "344×9=", 2) | Out-Null
$d.Content.Find.Execute("435×2=", $true, $false, $false, $false, $false, $true, 1, $false, "443×4=", 2) | Out-Null
$d.Content.Find.Execute("799×2=", $true, $false, $false, $false, $false, $true, 1, $false, "533×8=", 2) | Out-Null
$d.Content.Find.Execute("599×2=", $true, $false, $false, $false, $false, $true, 1, $false, "693×5=", 2) | Out-Null
$d.Content.Find.Execute("691×3=", $true, $false, $false, $false, $false, $true, 1, $false, "532×6=", 2) | Out-Null
$d.Content.Find.Execute("211×2=", $true, $false, $false, $false, $false, $true, 1, $false, "107×5=", 2) | Out-Null
$d.Content.Find.Execute("380×9=", $true, $false, $false, $false, $false, $true, 1, $false, "470×9=", 2) | Out-Null
$d.Content.Find.Execute("693×3=", $true, $false, $false, $false, $false, $true, 1, $false, "656×7=", 2) | Out-Null
$d.Content.Find.Execute("115×2=", $true, $false, $false, $false, $false, $true, 1, $false, "264×4=", 2) | Out-Null
$d.Content.Find.Execute("350×4=", $true, $false, $false, $false, $false, $true, 1, $false, "797×3=", 2) | Out-Null
$d.Content.Find.Execute("373×2=", $true, $false, $false, $false, $false, $true, 1, $false, "660×7=", 2) | Out-Null
$d.Content.Find.Execute("190×3=", $true, $false, $false, $false, $false, $true, 1, $false, "831×8=", 2) | Out-Null
$d.Content.Find.Execute("544×9=", $true, $false, $false, $false, $false, $true, 1, $false, "193×6=", 2) | Out-Null
$d.Content.Find.Execute("673×6=", $true, $false, $false, $false, $false, $true, 1, $false, "140×5=", 2) | Out-Null
$d.Content.Find.Execute("885×5=", $true, $false, $false, $false, $false, $true, 1, $false, "615×8=", 2) | Out-Null
$d.Content.Find.Execute("869×2=", $true, $false, $false, $false, $false, $true, 1, $false, "802×6=", 2) | Out-Null
$d.Content.Find.Execute("463×3=", $true, $false, $false, $false, $false, $true, 1, $false, "122×5=", 2) | Out-Null
$d.Content.Find.Execute("803×5=", $true, $false, $false, $false, $false, $true, 1, $false, "899×3=", 2) | Out-Null
$d.Content.Find.Execute("793×8=", $true, $false, $false, $false, $false, $true, 1, $false, "240×7=", 2) | Out-Null
$d.Content.Find.Execute("776×6=", $true, $false, $false, $false, $false, $true, 1, $false, "640×2=", 2) | Out-Null
